# Oneida slides (Oneida_slides_unicode.pptx) -- fill in the full
# Oneida-language phrases that were still placeholder single letters,
# and switch each touched run's Latin typeface from "NotoSans-Regular"
# to "Times New Roman" (per the target diff).
#
# Each run is addressed with TextRange.Characters(Start, Length) using
# the *original* run's Start/Length, which keeps the edit scoped to that
# single run (preserving bold/italic/underline/color/ea-typeface, etc.).
# Font is assigned before the text is replaced so the new (usually much
# longer) text inherits the updated formatting on that run instead of
# the replacement being split off into a fresh, default-formatted run.

$p = $ppt.ActivePresentation

function Set-RunText($range, $start, $length, $newText) {
    $c = $range.Characters($start, $length)
    $c.Font.Name = "Times New Roman"
    $c.Text = $newText
    return $newText.Length
}

# -----------------------------------------------------------------
# Slide 1
# -----------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Shape "TextShape 1" (sz=5400): "U" + "ha"
$sh = $s1.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pos = 1
$pos += Set-RunText $tr $pos 1 "Ukwehokú"
$pos += Set-RunText $tr $pos 2 "ha"

# Shape "TextShape 2" (sz=3200, italic): "T" "a" " " "y" " "
$sh = $s1.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pos = 1
$pos += Set-RunText $tr $pos 1 "Ta· "
$pos += Set-RunText $tr $pos 1 "aeswatahuhsiyósteʔ"
$pos += Set-RunText $tr $pos 1 " o·nʌ̀ tsiʔ náhteʔ ohʌ·tú "
$pos += Set-RunText $tr $pos 1 "yolihwatéhtuʔ"
$pos += Set-RunText $tr $pos 1 " "

# -----------------------------------------------------------------
# Slide 2
# -----------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Shape "TextShape 1" (sz=4800): "L" + "e"
$sh = $s2.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pos = 1
$pos += Set-RunText $tr $pos 1 "Lotihsóthaʔ Latihsakayu·téhs"
$pos += Set-RunText $tr $pos 1 "eʔ"

# Shape "TextShape 2" (sz=1600): 3 paragraphs, separated by a hard
# paragraph break (counts as one character, like \r, in TextRange.Text)
$sh = $s2.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pos = 1

# Paragraph 1: single run "T"
$pos += Set-RunText $tr $pos 1 "Ta· aeswélhekeʔ kʌtyókwaʔ né· tho niyo·lé· waʔkkwe·ní· né· kanehelatúksla."
$pos += 1  # paragraph break

# Paragraph 2: run "N", run "i" (underlined)
$pos += Set-RunText $tr $pos 1 "Né· katsaʔ ok nu· tukwató·ktʌ né· waʔtkatʌ̀·nukeʔ né· ʌskwatílhekeʔ né· tho niyo·lé waʔkkwe·ní· né· elhúwa wakewyʌtethaʔuhátyehseʔ né· kanʌ̀·laku akata·t"
$pos += Set-RunText $tr $pos 1 "í·"
$pos += 1  # paragraph break

# Paragraph 3: run "T", run "a?" (underlined, text unchanged), run "."
$pos += Set-RunText $tr $pos 1 "Ta· tho niyohtúhak né· yukwaʔnikúhl"
$pos += Set-RunText $tr $pos 2 "a?"
$pos += Set-RunText $tr $pos 1 ".  Táneʔ tho."

